$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("B2").Value = 0.1307692307692308
$ws.Range("C2").Value = 0.6538461538461539
$ws.Range("P2").Value = 0.09230769230769231
$ws.Range("S2").Value = 0.1230769230769231
$ws.Range("B3").Value = 0.005813953488372093
$ws.Range("C3").Value = 0.01744186046511628
$ws.Range("J3").Value = 0.02906976744186046
$ws.Range("P3").Value = 0.7441860465116279
$ws.Range("S3").Value = 0.2034883720930233
$ws.Range("J4").Value = 0.03636363636363636
$ws.Range("P4").Value = 0.7272727272727273
$ws.Range("S4").Value = 0.2363636363636364
$ws.Range("B6").Value = 0.03333333333333333
$ws.Range("D6").Value = 0.009523809523809525
$ws.Range("E6").Value = 0.004761904761904762
$ws.Range("F6").Value = 0.07142857142857142
$ws.Range("J6").Value = 0.2666666666666667
$ws.Range("O6").Value = 0.01428571428571429
$ws.Range("Q6").Value = 0.1666666666666667
$ws.Range("R6").Value = 0.09047619047619047
$ws.Range("S6").Value = 0.3428571428571429
$ws.Range("B7").Value = 0.1266666666666667
$ws.Range("D7").Value = 0.04666666666666667
$ws.Range("F7").Value = 0.02666666666666667
$ws.Range("J7").Value = 0.14
$ws.Range("O7").Value = 0.02
$ws.Range("Q7").Value = 0.1733333333333333
$ws.Range("R7").Value = 0.1066666666666667
$ws.Range("S7").Value = 0.36
$ws.Range("B8").Value = 0.1032110091743119
$ws.Range("D8").Value = 0.01834862385321101
$ws.Range("F8").Value = 0.07798165137614679
$ws.Range("J8").Value = 0.1055045871559633
$ws.Range("O8").Value = 0.02752293577981652
$ws.Range("Q8").Value = 0.1880733944954129
$ws.Range("R8").Value = 0.1215596330275229
$ws.Range("S8").Value = 0.3577981651376147
$ws.Range("B9").Value = 0.1055555555555556
$ws.Range("D9").Value = 0.02777777777777778
$ws.Range("E9").Value = 0.005555555555555556
$ws.Range("F9").Value = 0.04444444444444445
$ws.Range("J9").Value = 0.1111111111111111
$ws.Range("O9").Value = 0.03333333333333333
$ws.Range("Q9").Value = 0.1722222222222222
$ws.Range("R9").Value = 0.1055555555555556
$ws.Range("S9").Value = 0.3944444444444444
$ws.Range("B10").Value = 0.1071428571428571
$ws.Range("D10").Value = 0.02840909090909091
$ws.Range("E10").Value = 0.0008116883116883117
$ws.Range("F10").Value = 0.07386363636363637
$ws.Range("J10").Value = 0.1071428571428571
$ws.Range("O10").Value = 0.02191558441558442
$ws.Range("Q10").Value = 0.2021103896103896
$ws.Range("R10").Value = 0.1038961038961039
$ws.Range("S10").Value = 0.3547077922077922
$ws.Range("G11").Value = 0.1367521367521368
$ws.Range("J11").Value = 0.1111111111111111
$ws.Range("K11").Value = 0.1923076923076923
$ws.Range("L11").Value = 0.5555555555555556
$ws.Range("S11").Value = 0.004273504273504274
$ws.Range("G12").Value = 0.6617647058823529
$ws.Range("J12").Value = 0.2426470588235294
$ws.Range("L12").Value = 0.03676470588235294
$ws.Range("S12").Value = 0.05882352941176471
$ws.Range("G13").Value = 0.6976744186046512
$ws.Range("J13").Value = 0.2558139534883721
$ws.Range("S13").Value = 0.04651162790697674
$ws.Range("F15").Value = 0.03347280334728033
$ws.Range("H15").Value = 0.2050209205020921
$ws.Range("I15").Value = 0.06694560669456066
$ws.Range("J15").Value = 0.3221757322175732
$ws.Range("K15").Value = 0.05857740585774059
$ws.Range("M15").Value = 0.004184100418410041
$ws.Range("O15").Value = 0.04602510460251046
$ws.Range("S15").Value = 0.2635983263598327
$ws.Range("F16").Value = 0.01063829787234043
$ws.Range("H16").Value = 0.1808510638297872
$ws.Range("I16").Value = 0.1117021276595745
$ws.Range("J16").Value = 0.3882978723404255
$ws.Range("K16").Value = 0.0797872340425532
$ws.Range("M16").Value = 0.01595744680851064
$ws.Range("O16").Value = 0.05851063829787234
$ws.Range("S16").Value = 0.1542553191489362
$ws.Range("F17").Value = 0.02137767220902613
$ws.Range("H17").Value = 0.166270783847981
$ws.Range("I17").Value = 0.1021377672209026
$ws.Range("J17").Value = 0.4513064133016627
$ws.Range("K17").Value = 0.07125890736342043
$ws.Range("M17").Value = 0.02612826603325416
$ws.Range("O17").Value = 0.04513064133016627
$ws.Range("S17").Value = 0.1163895486935867
$ws.Range("F18").Value = 0.03418803418803419
$ws.Range("H18").Value = 0.1666666666666667
$ws.Range("I18").Value = 0.06837606837606838
$ws.Range("J18").Value = 0.4786324786324787
$ws.Range("K18").Value = 0.0811965811965812
$ws.Range("M18").Value = 0.008547008547008548
$ws.Range("O18").Value = 0.0641025641025641
$ws.Range("S18").Value = 0.09829059829059829
$ws.Range("F19").Value = 0.01380500431406385
$ws.Range("H19").Value = 0.2131147540983606
$ws.Range("I19").Value = 0.07420189818809318
$ws.Range("J19").Value = 0.3805004314063848
$ws.Range("K19").Value = 0.09663503019844694
$ws.Range("M19").Value = 0.02415875754961173
$ws.Range("N19").Value = 0.001725625539257981
$ws.Range("O19").Value = 0.08369283865401209
$ws.Range("S19").Value = 0.1121656600517688
